$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "ежедневный посыл" trigger time in row 14 (column B)
$ws.Range("B14").Value = "23:49 - 23:59"

# Move the active selection to B22, matching the saved view state
$ws.Range("B22").Select()
